$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation for "Haba" was recorded for Femacal de
# La Calera. Insert a fresh row above the current row 250 (pushing every
# row from 250 down through 311 to 251 through 312, A1:R311 -> A1:R312)
# and populate it with the new entry.
$ws.Rows(250).Insert()

$ws.Cells.Item(250, 1).Value  = 3
$ws.Cells.Item(250, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(250, 3).Value  = "Coquimbo"
$ws.Cells.Item(250, 4).Value  = 45204
$ws.Cells.Item(250, 5).Value  = 5
$ws.Cells.Item(250, 6).Value  = 100112026
$ws.Cells.Item(250, 7).Value  = "Haba"
$ws.Cells.Item(250, 8).Value  = "Sin especificar"
$ws.Cells.Item(250, 9).Value  = "Primera"
$ws.Cells.Item(250, 10).Value = 80
$ws.Cells.Item(250, 11).Value = 13000
$ws.Cells.Item(250, 12).Value = 13500
$ws.Cells.Item(250, 13).Value = 13250
$ws.Cells.Item(250, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(250, 15).Value = "Provincia de Limar" + [char]0x00ED
$ws.Cells.Item(250, 16).Value = 530
$ws.Cells.Item(250, 17).Value = 25
$ws.Cells.Item(250, 18).Value = "Hortaliza"
